# Realestate Update resale numbers 2023-06-25 12:34
# Append a new data row (row 76) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

# Columns A-D hold text in the source data (date/time/weekday/week-number
# stored as strings, not native date/number types). Force text formatting
# before assignment so Excel doesn't auto-coerce them to a date serial or
# a number, then clear the formatting so no stray style is left behind on
# the new cells (matches the plain, unstyled cells used by existing rows).
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Range("A$row").Value = "2023-06-25"
$ws.Range("B$row").Value = "12:34:05"
$ws.Range("C$row").Value = "Sunday"
$ws.Range("D$row").Value = "26"

$textRange.ClearFormats()

# Columns E-T are plain numbers.
$ws.Range("E$row").Value = 122644
$ws.Range("F$row").Value = 134264
$ws.Range("G$row").Value = 163115
$ws.Range("H$row").Value = 133540
$ws.Range("I$row").Value = 177415
$ws.Range("J$row").Value = 116438
$ws.Range("K$row").Value = 203200
$ws.Range("L$row").Value = 225861
$ws.Range("M$row").Value = 175647
$ws.Range("N$row").Value = 104228
$ws.Range("O$row").Value = 39545
$ws.Range("P$row").Value = 33810
$ws.Range("Q$row").Value = 52050
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 35850
$ws.Range("T$row").Value = -1
